# Remove the second slide ("My new slide") from the presentation.
# This corresponds to the slide with sldId 258 (r:id R0e34d47079a84dcb),
# which currently sits between slide1 and slide2 in the slide order.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$s.Delete()
